$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add REGIAO (C1) and CODREG (D1)
$ws.Range("C1").Value = "REGIAO"
$ws.Range("D1").Value = "CODREG"

# Mapping of row -> (Region name, Region code)
$regiaoData = @{
    2  = @("Norte", 1)
    3  = @("Norte", 1)
    4  = @("Norte", 1)
    5  = @("Norte", 1)
    6  = @("Norte", 1)
    7  = @("Norte", 1)
    8  = @("Norte", 1)
    9  = @("Nordeste", 2)
    10 = @("Nordeste", 2)
    11 = @("Nordeste", 2)
    12 = @("Nordeste", 2)
    13 = @("Nordeste", 2)
    14 = @("Nordeste", 2)
    15 = @("Nordeste", 2)
    16 = @("Nordeste", 2)
    17 = @("Nordeste", 2)
    18 = @("Sudeste", 3)
    19 = @("Sudeste", 3)
    20 = @("Sudeste", 3)
    21 = @("Sudeste", 3)
    22 = @("Sul", 4)
    23 = @("Sul", 4)
    24 = @("Sul", 4)
    25 = @("Centro-oeste", 5)
    26 = @("Centro-oeste", 5)
    27 = @("Centro-oeste", 5)
    28 = @("Centro-oeste", 5)
}

foreach ($row in 2..28) {
    $pair = $regiaoData[$row]
    $ws.Cells.Item($row, 3).Value = $pair[0]
    $ws.Cells.Item($row, 4).Value = $pair[1]
}

# Match author's formatting pass: auto-size the new REGIAO column to fit
# its longest entry ("Centro-oeste"), same as cols A/B were already sized.
$ws.Columns.Item(3).AutoFit()
